$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Map of cell address -> new value, reflecting updated registration counts.
$updates = @{
    "E3"  = 27;  "F3"  = 24;  "H3"  = 24
    "E4"  = 31
    "E5"  = 92;  "F5"  = 58;  "H5"  = 58
    "E6"  = 35
    "E10" = 306; "F10" = 152; "H10" = 152
    "E11" = 213
    "E12" = 315
    "E13" = 94
    "E14" = 88
    "E18" = 43
    "E20" = 68
    "E21" = 102
    "E22" = 120
    "E23" = 125
    "E24" = 139
    "E25" = 159
    "E26" = 88
    "E27" = 208; "F27" = 107; "H27" = 107
    "E28" = 122
    "E29" = 123
    "E30" = 144
    "E31" = 53
    "E32" = 128
    "E33" = 196
    "E34" = 147
    "E35" = 95
    "E37" = 106
    "E39" = 132
    "E40" = 181
    "E41" = 253; "F41" = 118; "H41" = 118
    "E42" = 233
    "E43" = 74
    "E44" = 211
    "E45" = 86
    "E46" = 199; "F46" = 112; "H46" = 112
    "E47" = 297; "F47" = 151; "H47" = 151
    "E48" = 139
    "E49" = 180
    "E50" = 155
    "E51" = 143
    "E52" = 18
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
